$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column T (20) mirrors column S (19) formatting for each data row.
# Copy the style of each existing S-column cell into the new T-column cell,
# then set the 2023 data value for that row.
$tData = @(
    @{ Row = 3; Value = 2023 },
    @{ Row = 4; Value = 2.3381104968484805 },
    @{ Row = 5; Value = 2.0344672190198714 },
    @{ Row = 6; Value = 2.6483752218014245 },
    @{ Row = 7; Value = 3.9852372948902328 },
    @{ Row = 8; Value = 4.5532396299967433 },
    @{ Row = 9; Value = 3.4291318466903733 },
    @{ Row = 10; Value = 1.2089851778417198 },
    @{ Row = 11; Value = 1.521116134174612 },
    @{ Row = 12; Value = 0.9008846687447073 },
    @{ Row = 13; Value = 3.694303753043183 },
    @{ Row = 14; Value = 4.0607488020791038 },
    @{ Row = 15; Value = 3.327319511401615 },
    @{ Row = 16; Value = 0.32236434908190637 },
    @{ Row = 17; Value = 0 },
    @{ Row = 18; Value = 0.63756806039044667 },
    @{ Row = 19; Value = 2.1691385808410835 },
    @{ Row = 20; Value = 1.5024572004578396 },
    @{ Row = 21; Value = 2.8259763748375066 },
    @{ Row = 22; Value = 6.1744985943935555 },
    @{ Row = 23; Value = 4.3993752887090034 },
    @{ Row = 24; Value = 7.9169155696940479 },
    @{ Row = 25; Value = 2.8763040791558883 },
    @{ Row = 26; Value = 1.4751329463567904 },
    @{ Row = 27; Value = 4.2954684675262591 },
    @{ Row = 28; Value = 1.8177568880002077 },
    @{ Row = 29; Value = 1.581380197008345 },
    @{ Row = 30; Value = 2.103608453446189 },
    @{ Row = 31; Value = 1.3736037318066185 },
    @{ Row = 32; Value = 2.249820014398848 },
    @{ Row = 33; Value = 0.53701655085009725 }
)

foreach ($item in $tData) {
    $r = $item.Row
    $sCell = $ws.Cells.Item($r, 19)
    $tCell = $ws.Cells.Item($r, 20)

    # Copy number format / font / border from the column-S cell in this row
    $sCell.Copy()
    $tCell.PasteSpecial(-4122)

    $tCell.Value = $item.Value
}

# Row 4 grew a hair taller once the 2023 column was populated.
$ws.Rows.Item(4).RowHeight = 16.5

# Clear the clipboard marquee left behind by the Copy() calls above.
$excel.CutCopyMode = $false

